$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.022.60'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.869.71'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.93'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5156'
$ws.Range("E7").Value = '  +2.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3849'
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08278'
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.112'
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.51'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.207'
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.55'
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.859.88'
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.298'
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001100'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.88'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06640'
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.030'
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").Value = '28.068.89'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.10'
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.250'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").Value = '2.082.33'
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.500'
$ws.Range("E27").Value = '  -3.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.43'
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.56'
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.13'
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1068'
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.817'
$ws.Range("E33").Value = '  +3.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.594'
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.482'
$ws.Range("E35").Value = '  -1.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02426'
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06508'
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("E38").Value = '  +0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6582'
$ws.Range("E39").Value = '  +2.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.201'
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("E41").Value = '  +2.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.208'
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.20'
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("E44").Value = '  +2.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.08'
$ws.Range("E45").Value = '  -1.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.285'
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.670'
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.220'
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.93'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.54'
